$wb = $excel.ActiveWorkbook

$daily = $wb.Worksheets.Item("Daily")
$hourly = $wb.Worksheets.Item("Hourly")

# Daily sheet - row 2 updates
$daily.Range("G2").Value = 2989.75
$daily.Range("H2").Value = 6264.64
$daily.Range("I2").Value = 738.13
$daily.Range("J2").Value = 2981.27
$daily.Range("K2").Value = 5900.62
$daily.Range("L2").Value = 739.11

# Hourly sheet - rows 9-19 updates
$hourly.Range("I9").Value = 62.87
$hourly.Range("K9").Value = 8.34
$hourly.Range("L9").Value = 15.95
$hourly.Range("M9").Value = 7.39
$hourly.Range("H10").Value = 112.46
$hourly.Range("I10").Value = 441.61
$hourly.Range("J10").Value = 49.96
$hourly.Range("K10").Value = 111.86
$hourly.Range("L10").Value = 397.98
$hourly.Range("M10").Value = 45.94
$hourly.Range("H11").Value = 254.42
$hourly.Range("I11").Value = 639.5700000000001
$hourly.Range("J11").Value = 73.03
$hourly.Range("K11").Value = 254.06
$hourly.Range("L11").Value = 616.83
$hourly.Range("M11").Value = 71.31999999999999
$hourly.Range("H12").Value = 375.25
$hourly.Range("I12").Value = 735.02
$hourly.Range("J12").Value = 86.44
$hourly.Range("K12").Value = 374.87
$hourly.Range("L12").Value = 717.77
$hourly.Range("M12").Value = 85.73999999999999
$hourly.Range("H13").Value = 456.2
$hourly.Range("I13").Value = 782.58
$hourly.Range("J13").Value = 93.84999999999999
$hourly.Range("K13").Value = 455.85
$hourly.Range("L13").Value = 759.33
$hourly.Range("M13").Value = 97.43000000000001
$hourly.Range("H14").Value = 487.57
$hourly.Range("I14").Value = 798.66
$hourly.Range("J14").Value = 96.48
$hourly.Range("K14").Value = 487.28
$hourly.Range("L14").Value = 772.16
$hourly.Range("M14").Value = 102.38
$hourly.Range("H15").Value = 465.89
$hourly.Range("I15").Value = 787.73
$hourly.Range("J15").Value = 94.66
$hourly.Range("K15").Value = 465.63
$hourly.Range("L15").Value = 763.72
$hourly.Range("M15").Value = 98.90000000000001
$hourly.Range("H16").Value = 393.54
$hourly.Range("I16").Value = 746.77
$hourly.Range("J16").Value = 88.18000000000001
$hourly.Range("K16").Value = 392.43
$hourly.Range("L16").Value = 725.37
$hourly.Range("M16").Value = 88.76000000000001
$hourly.Range("H17").Value = 279.02
$hourly.Range("I17").Value = 662.4400000000001
$hourly.Range("J17").Value = 76.04000000000001
$hourly.Range("K17").Value = 275.97
$hourly.Range("L17").Value = 626.3099999999999
$hourly.Range("M17").Value = 76.53
$hourly.Range("H18").Value = 138.85
$hourly.Range("I18").Value = 492
$hourly.Range("J18").Value = 55.32
$hourly.Range("K18").Value = 137.13
$hourly.Range("L18").Value = 449.9
$hourly.Range("M18").Value = 51.78
$hourly.Range("I19").Value = 115.39
$hourly.Range("K19").Value = 17.85
$hourly.Range("L19").Value = 55.3
$hourly.Range("M19").Value = 12.95
